# Refresh the crypto price/volume table with the latest scrape results.
# (Row 42-48 in the source also reshuffled ranking order; B/C/D/E are
# rewritten per-row to match, including two rows that simply swapped places.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.887.89'  # D2
$ws.Cells.Item(2, 5).Value = '  +5.25%  '  # E2

$ws.Cells.Item(3, 4).Value = '3.517.36'  # D3
$ws.Cells.Item(3, 5).Value = '  +3.43%  '  # E3

$ws.Cells.Item(4, 5).Value = '  -0.04%  '  # E4

$ws.Cells.Item(5, 4).Value = '''592.57'  # D5
$ws.Cells.Item(5, 5).Value = '  +4.42%  '  # E5

$ws.Cells.Item(6, 4).Value = '''169.98'  # D6
$ws.Cells.Item(6, 5).Value = '  +8.72%  '  # E6

$ws.Cells.Item(8, 4).Value = '3.520.72'  # D8
$ws.Cells.Item(8, 5).Value = '  +3.48%  '  # E8

$ws.Cells.Item(9, 4).Value = '''0.579'  # D9
$ws.Cells.Item(9, 5).Value = '  +2.21%  '  # E9

$ws.Cells.Item(10, 4).Value = '''7.25'  # D10
$ws.Cells.Item(10, 5).Value = '  +0.12%  '  # E10

$ws.Cells.Item(11, 4).Value = '''0.125'  # D11
$ws.Cells.Item(11, 5).Value = '  +5.32%  '  # E11

$ws.Cells.Item(12, 4).Value = '''0.437'  # D12
$ws.Cells.Item(12, 5).Value = '  +3.60%  '  # E12

$ws.Cells.Item(13, 4).Value = '4.113.78'  # D13
$ws.Cells.Item(13, 5).Value = '  +3.18%  '  # E13

$ws.Cells.Item(14, 5).Value = '  +0.46%  '  # E14

$ws.Cells.Item(15, 4).Value = '''28.07'  # D15
$ws.Cells.Item(15, 5).Value = '  +4.51%  '  # E15

$ws.Cells.Item(16, 4).Value = '66.833.96'  # D16
$ws.Cells.Item(16, 5).Value = '  +5.00%  '  # E16

$ws.Cells.Item(17, 4).Value = '''0.0000178'  # D17
$ws.Cells.Item(17, 5).Value = '  +4.50%  '  # E17

$ws.Cells.Item(18, 4).Value = '3.493.05'  # D18
$ws.Cells.Item(18, 5).Value = '  +2.56%  '  # E18

$ws.Cells.Item(19, 4).Value = '''6.29'  # D19
$ws.Cells.Item(19, 5).Value = '  +3.62%  '  # E19

$ws.Cells.Item(20, 4).Value = '''14.01'  # D20
$ws.Cells.Item(20, 5).Value = '  +3.94%  '  # E20

$ws.Cells.Item(21, 4).Value = '''390.13'  # D21
$ws.Cells.Item(21, 5).Value = '  +2.04%  '  # E21

$ws.Cells.Item(22, 5).Value = '  +2.75%  '  # E22

$ws.Cells.Item(23, 4).Value = '''73.21'  # D23
$ws.Cells.Item(23, 5).Value = '  +2.96%  '  # E23

$ws.Cells.Item(24, 5).Value = '  +0.53%  '  # E24

$ws.Cells.Item(25, 4).Value = '''0.528'  # D25
$ws.Cells.Item(25, 5).Value = '  +2.91%  '  # E25

$ws.Cells.Item(26, 4).Value = '''0.0000124'  # D26
$ws.Cells.Item(26, 5).Value = '  +9.08%  '  # E26

$ws.Cells.Item(27, 4).Value = '''10.13'  # D27
$ws.Cells.Item(27, 5).Value = '  +4.64%  '  # E27

$ws.Cells.Item(28, 5).Value = '  +1.87%  '  # E28

$ws.Cells.Item(29, 4).Value = '''0.998'  # D29
$ws.Cells.Item(29, 5).Value = '  -0.14%  '  # E29

$ws.Cells.Item(30, 4).Value = '''6.39'  # D30
$ws.Cells.Item(30, 5).Value = '  +6.20%  '  # E30

$ws.Cells.Item(31, 4).Value = '''1.47'  # D31
$ws.Cells.Item(31, 5).Value = '  +5.83%  '  # E31

$ws.Cells.Item(32, 5).Value = '  +4.41%  '  # E32

$ws.Cells.Item(33, 4).Value = '''23.52'  # D33
$ws.Cells.Item(33, 5).Value = '  +3.11%  '  # E33

$ws.Cells.Item(34, 4).Value = '''7.41'  # D34
$ws.Cells.Item(34, 5).Value = '  +7.49%  '  # E34

$ws.Cells.Item(35, 4).Value = '''0.999'  # D35
$ws.Cells.Item(35, 5).Value = '  +0.06%  '  # E35

$ws.Cells.Item(36, 4).Value = '''1.58'  # D36
$ws.Cells.Item(36, 5).Value = '  +5.29%  '  # E36

$ws.Cells.Item(37, 4).Value = '''161.19'  # D37
$ws.Cells.Item(37, 5).Value = '  +0.48%  '  # E37

$ws.Cells.Item(38, 4).Value = '''0.908'  # D38
$ws.Cells.Item(38, 5).Value = '  +8.04%  '  # E38

$ws.Cells.Item(39, 5).Value = '  +6.70%  '  # E39

$ws.Cells.Item(40, 4).Value = '''0.0747'  # D40
$ws.Cells.Item(40, 5).Value = '  +4.59%  '  # E40

$ws.Cells.Item(41, 4).Value = '''4.66'  # D41
$ws.Cells.Item(41, 5).Value = '  +7.61%  '  # E41

$ws.Cells.Item(42, 2).Value = 'EnergySwap'  # B42
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'  # C42
$ws.Cells.Item(42, 4).Value = '''26.54'  # D42
$ws.Cells.Item(42, 5).Value = '  +2.63%  '  # E42

$ws.Cells.Item(43, 2).Value = 'RenderToken'  # B43
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'  # C43
$ws.Cells.Item(43, 4).Value = '''6.72'  # D43
$ws.Cells.Item(43, 5).Value = '  +5.56%  '  # E43

$ws.Cells.Item(44, 2).Value = 'Maker'  # B44
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'  # C44
$ws.Cells.Item(44, 4).Value = '2.807.11'  # D44
$ws.Cells.Item(44, 5).Value = '  +0.01%  '  # E44

$ws.Cells.Item(45, 2).Value = 'OKB'  # B45
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'  # C45
$ws.Cells.Item(45, 4).Value = '''43.51'  # D45
$ws.Cells.Item(45, 5).Value = '  +1.36%  '  # E45

$ws.Cells.Item(46, 4).Value = '''26.83'  # D46
$ws.Cells.Item(46, 5).Value = '  +5.00%  '  # E46

$ws.Cells.Item(47, 2).Value = 'dogwifhat'  # B47
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'  # C47
$ws.Cells.Item(47, 4).Value = '''2.55'  # D47
$ws.Cells.Item(47, 5).Value = '  +10.44%  '  # E47

$ws.Cells.Item(48, 2).Value = 'VeChain'  # B48
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'  # C48
$ws.Cells.Item(48, 4).Value = '''0.0313'  # D48
$ws.Cells.Item(48, 5).Value = '  +3.66%  '  # E48

$ws.Cells.Item(49, 4).Value = '''355.31'  # D49
$ws.Cells.Item(49, 5).Value = '  +8.83%  '  # E49

$ws.Cells.Item(50, 5).Value = '  +6.51%  '  # E50

$ws.Cells.Item(51, 4).Value = '''33.20'  # D51
$ws.Cells.Item(51, 5).Value = '  +10.76%  '  # E51
